$d = $word.ActiveDocument

# The last paragraph in the document contains only the _GoBack bookmark.
# We need to add new formatted text runs before that bookmark (within the
# same paragraph) and then append a brand-new empty paragraph after it.
$p = $d.Paragraphs.Last
$r = $p.Range

$xml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t>Credit information, debit finance,</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>376254468525489</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>, credit finance</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>, card information</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$r.InsertXML($xml)
